# "New Microsite scripts support to Beta server"
#
# AMSIN sheet: append 4 new registration rows (13-16) after the existing
# row 12, matching its row-level styling (s="5" text/number cells,
# s="10" datetime cell for column B).
#
# AMS sheet: the existing row 12 ("payy166") picks up explicit s="5"
# styling (it previously had none) and its run-time (column B) value is
# refreshed; a brand-new row 13 ("betapay167") is appended below it with
# no explicit styling (matching the row 12 had before this edit).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a literal (non-date-coerced) text value into a cell while
# preserving whatever style is already sitting on that cell. Plain
# `Range.Value = "2022-09-15"` gets auto-recognised as a date by the
# engine (same smart-parsing Excel itself does) which both changes the
# stored value AND mints a brand-new number-format style. Instead, stash
# the literal text as a formula result (`="2022-09-15"`, which is never
# re-interpreted as a date) in a throwaway scratch cell, then copy only
# the computed VALUE (xlPasteValues = -4163) onto the destination — this
# leaves the destination's existing style untouched and doesn't leave any
# stray number-format behind in styles.xml either.
# ---------------------------------------------------------------------
function Set-LiteralText {
    param($ws, $row, $col, $text)

    $scratchRow = 500
    $scratch = $ws.Cells.Item($scratchRow, $col)
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
    $ws.Rows.Item($scratchRow).Delete()
}

function Set-DataRow {
    param($ws, $row, $date, $time, $name, $total, $pass, $fail, $taken)

    Set-LiteralText $ws $row 1 $date
    $ws.Cells.Item($row, 2).Value = $time
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $taken
}

# ---------------------------------------------------------------------
# AMSIN sheet — append rows 13..16
# ---------------------------------------------------------------------
$amsin = $wb.Worksheets.Item("AMSIN")

# Duplicate row 12 (style template: s="5" / s="10") four times, inserting
# each copy immediately below the existing data so every new row inherits
# the exact same per-cell styling as row 12.
$amsin.Rows.Item(12).Copy()
$amsin.Rows.Item(13).Insert()
$amsin.Rows.Item(12).Copy()
$amsin.Rows.Item(14).Insert()
$amsin.Rows.Item(12).Copy()
$amsin.Rows.Item(15).Insert()
$amsin.Rows.Item(12).Copy()
$amsin.Rows.Item(16).Insert()

Set-DataRow $amsin 13 "2022-09-15" 44819.63156479166 "payecs166"   41 41 0 1.44
Set-DataRow $amsin 14 "2022-09-16" 44820.64721446759 "fstcpay167"  41 41 0 1.55
Set-DataRow $amsin 15 "2022-09-19" 44823.60005887732 "scndcpay167" 41 40 1 2.66
Set-DataRow $amsin 16 "2022-09-20" 44824.38984690973 "finalpay167" 41 41 0 1.16

# ---------------------------------------------------------------------
# AMS sheet — restyle row 12, refresh its run-time, append row 13
# ---------------------------------------------------------------------
$ams = $wb.Worksheets.Item("AMS")

# Row 11 already carries the s="5"/s="10" styling row 12 needs to gain.
# Copy-inserting it at position 12 shifts the current (unstyled) row 12
# down to row 13 untouched, and stamps the new row 12 with row 11's
# per-cell styles.
$ams.Rows.Item(11).Copy()
$ams.Rows.Item(12).Insert()

Set-DataRow $ams 12 "2022-09-08" 44812.55046903935 "payy166" 41 41 0 1.18

# Row 13 (the original row-12 data, shifted down by the insert above)
# keeps its old bare styling; overwrite it with the new Beta-server
# registration entry.
Set-DataRow $ams 13 "2022-09-20" 44824.73223812653 "betapay167" 41 41 0 1
